# KIBON-2653: Change Column kinderFaktor LATS Statistik
#
# The "Tagesschulen" sheet had a single "2_Kinder_Faktor" column
# (placeholder {kinderFaktor}). It is replaced by two columns:
# "2_Kinder_Faktor_3" / {kinderFaktor3} and "2_Kinder_Faktor_15" / {kinderFaktor15}.

$wb = $excel.ActiveWorkbook
$wsGemeinden = $wb.Worksheets.Item("Gemeinden")
$wsTages = $wb.Worksheets.Item("Tagesschulen")

# Insert a new column at I so the old "2_Kinder_Faktor" column (I) shifts to J,
# leaving a fresh empty column I to hold the new "3er" factor column.
$wsTages.Columns.Item(9).Insert()

# Fill in the data row first, then the header row, so new shared-string
# entries get appended to the table in the same order Excel produced them.
$wsTages.Range("I2").Value = "{kinderFaktor3}"
$wsTages.Range("J2").Value = "{kinderFaktor15}"
$wsTages.Range("I1").Value = "2_Kinder_Faktor_3"
$wsTages.Range("J1").Value = "2_Kinder_Faktor_15"

# Restore cursor/selection state: Gemeinden keeps its last selection (O9),
# Tagesschulen becomes the active (visible) sheet with J2 selected.
$wsGemeinden.Activate() | Out-Null
$wsGemeinden.Range("O9").Select() | Out-Null
$wsTages.Activate() | Out-Null
$wsTages.Range("J2").Select() | Out-Null
